# Daily attendance processing - normalize "Recorded By" (column G) values
# so that automated/system accounts are listed before the human reviewer
# email address in the comma-separated attribution list.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Exact literal normalizations observed for the "Recorded By" column.
# Key = current full cell text, Value = normalized replacement text.
$map = @{
    "backup@backdoor.com, System"         = "System, backup@backdoor.com"
    "backup@backdoor.com, System, system" = "system, System, backup@backdoor.com"
    "dnasr281@gmail.com, System"          = "System, dnasr281@gmail.com"
}

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)  # Column G = "Recorded By"
    $current = $cell.Text

    if ($map.ContainsKey($current)) {
        $cell.Value = $map[$current]
    }
}
